# Revert "Powerpoint writer: consolidate text run nodes."
#
# Split the title runs on slide 1 ("Example " + "numbering " + "MWE") and
# slide 2 ("A " + "second " + "slide") so the trailing space of each of the
# first two words becomes its own run, instead of being merged into the
# word that precedes it. Re-assigning the same text to a `Characters`
# sub-range forces the writer to break the run at that boundary without
# changing the visible text.

$p = $ppt.ActivePresentation

# Slide 1 title: "Example " | "numbering " | "MWE"
#            -> "Example" | " " | "numbering" | " " | "MWE"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(8, 1).Text = " "
$tr1.Characters(18, 1).Text = " "

# Slide 2 title: "A " | "second " | "slide"
#            -> "A" | " " | "second" | " " | "slide"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(2, 1).Text = " "
$tr2.Characters(9, 1).Text = " "
